$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 3) down to the new row 4
$ws.Range("A3:H3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new trade data for the 20 minute trade
$ws.Range("A4").Value = 9807.8799999999992
$ws.Range("B4").Value = 9881
$ws.Range("C4").Value = 20.3
$ws.Range("D4").Value = 20.149999999999999
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = -0.74
$ws.Range("G4").Value = 42608.640486111108
$ws.Range("H4").Value = $false
